$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows 66/67: Q column 1 -> 0 ---
$ws.Range("Q66").Value = 0
$ws.Range("Q67").Value = 0

# --- Fix existing rows 1142/1143: R column inlineStr("") -> numeric 0 ---
$ws.Range("R1142").Value = 0
$ws.Range("R1143").Value = 0

# --- Capture the date number format used by column A (style index 2) ---
$dateFmt = $ws.Range("A1143").NumberFormat()

# --- Append new weekly rows 1144-1169 ---
# Row 1144
$ws.Range("A1144").Value = 45474
$ws.Range("B1144").Value = 314.3233939346234
$ws.Range("C1144").Value = 345.2095218014308
$ws.Range("D1144").Value = 313.3799150682936
$ws.Range("E1144").Value = 333.4410095214844
$ws.Range("F1144").Value = "'"
$ws.Range("F1144").Style = "Normal"
$ws.Range("G1144").Value = 46291955
$ws.Range("H1144").Value = 2024
$ws.Range("I1144").Value = 7
$ws.Range("J1144").Value = 1
$ws.Range("K1144").Value = 0
$ws.Range("L1144").Value = 0
$ws.Range("M1144").Value = 0
$ws.Range("N1144").Value = 27
$ws.Range("O1144").Value = 0
$ws.Range("P1144").Value = 0
$ws.Range("Q1144").Value = 0
$ws.Range("R1144").Value = "'"
$ws.Range("R1144").Style = "Normal"
$ws.Range("A1144").NumberFormat = $dateFmt

# Row 1145
$ws.Range("A1145").Value = 45481
$ws.Range("B1145").Value = 335.6755231471197
$ws.Range("C1145").Value = 346.0536690479748
$ws.Range("D1145").Value = 317.2034320839919
$ws.Range("E1145").Value = 324.4532470703125
$ws.Range("F1145").Value = "'"
$ws.Range("F1145").Style = "Normal"
$ws.Range("G1145").Value = 18963592
$ws.Range("H1145").Value = 2024
$ws.Range("I1145").Value = 7
$ws.Range("J1145").Value = 8
$ws.Range("K1145").Value = 0
$ws.Range("L1145").Value = 0
$ws.Range("M1145").Value = 0
$ws.Range("N1145").Value = 28
$ws.Range("O1145").Value = 0
$ws.Range("P1145").Value = 0
$ws.Range("Q1145").Value = 0
$ws.Range("R1145").Value = "'"
$ws.Range("R1145").Style = "Normal"
$ws.Range("A1145").NumberFormat = $dateFmt

# Row 1146
$ws.Range("A1146").Value = 45488
$ws.Range("B1146").Value = 325.14842673455
$ws.Range("C1146").Value = 334.2354950776646
$ws.Range("D1146").Value = 305.8818488124497
$ws.Range("E1146").Value = 307.0736083984375
$ws.Range("F1146").Value = "'"
$ws.Range("F1146").Style = "Normal"
$ws.Range("G1146").Value = 18050304
$ws.Range("H1146").Value = 2024
$ws.Range("I1146").Value = 7
$ws.Range("J1146").Value = 15
$ws.Range("K1146").Value = 0
$ws.Range("L1146").Value = 0
$ws.Range("M1146").Value = 0
$ws.Range("N1146").Value = 29
$ws.Range("O1146").Value = 0
$ws.Range("P1146").Value = 0
$ws.Range("Q1146").Value = 0
$ws.Range("R1146").Value = "'"
$ws.Range("R1146").Style = "Normal"
$ws.Range("A1146").NumberFormat = $dateFmt

# Row 1147
$ws.Range("A1147").Value = 45495
$ws.Range("B1147").Value = 307.0736061111193
$ws.Range("C1147").Value = 347.2950826845847
$ws.Range("D1147").Value = 293.467810554544
$ws.Range("E1147").Value = 335.2782897949219
$ws.Range("F1147").Value = "'"
$ws.Range("F1147").Style = "Normal"
$ws.Range("G1147").Value = 79505331
$ws.Range("H1147").Value = 2024
$ws.Range("I1147").Value = 7
$ws.Range("J1147").Value = 22
$ws.Range("K1147").Value = 0
$ws.Range("L1147").Value = 0
$ws.Range("M1147").Value = 0
$ws.Range("N1147").Value = 30
$ws.Range("O1147").Value = 2
$ws.Range("P1147").Value = 0
$ws.Range("Q1147").Value = 0
$ws.Range("R1147").Value = "'"
$ws.Range("R1147").Style = "Normal"
$ws.Range("A1147").NumberFormat = $dateFmt

# Row 1148
$ws.Range("A1148").Value = 45502
$ws.Range("B1148").Value = 335.3279349121502
$ws.Range("C1148").Value = 361.9932903447666
$ws.Range("D1148").Value = 331.7526864768337
$ws.Range("E1148").Value = 341.7832336425781
$ws.Range("F1148").Value = "'"
$ws.Range("F1148").Style = "Normal"
$ws.Range("G1148").Value = 49114124
$ws.Range("H1148").Value = 2024
$ws.Range("I1148").Value = 7
$ws.Range("J1148").Value = 29
$ws.Range("K1148").Value = 0
$ws.Range("L1148").Value = 0
$ws.Range("M1148").Value = 0
$ws.Range("N1148").Value = 31
$ws.Range("O1148").Value = 1
$ws.Range("P1148").Value = 0
$ws.Range("Q1148").Value = 0
$ws.Range("R1148").Value = "'"
$ws.Range("R1148").Style = "Normal"
$ws.Range("A1148").NumberFormat = $dateFmt

# Row 1149
$ws.Range("A1149").Value = 45509
$ws.Range("B1149").Value = 323.0628677171563
$ws.Range("C1149").Value = 332.3485729072133
$ws.Range("D1149").Value = 307.1729012940274
$ws.Range("E1149").Value = 321.7718200683594
$ws.Range("F1149").Value = "'"
$ws.Range("F1149").Style = "Normal"
$ws.Range("G1149").Value = 42528703
$ws.Range("H1149").Value = 2024
$ws.Range("I1149").Value = 8
$ws.Range("J1149").Value = 5
$ws.Range("K1149").Value = 0
$ws.Range("L1149").Value = 0
$ws.Range("M1149").Value = 0
$ws.Range("N1149").Value = 32
$ws.Range("O1149").Value = 0
$ws.Range("P1149").Value = 0
$ws.Range("Q1149").Value = 0
$ws.Range("R1149").Value = "'"
$ws.Range("R1149").Style = "Normal"
$ws.Range("A1149").NumberFormat = $dateFmt

# Row 1150
$ws.Range("A1150").Value = 45516
$ws.Range("B1150").Value = 320.9773121078645
$ws.Range("C1150").Value = 325.049091481558
$ws.Range("D1150").Value = 306.3783943233465
$ws.Range("E1150").Value = 307.7191162109375
$ws.Range("F1150").Value = "'"
$ws.Range("F1150").Style = "Normal"
$ws.Range("G1150").Value = 19009834
$ws.Range("H1150").Value = 2024
$ws.Range("I1150").Value = 8
$ws.Range("J1150").Value = 12
$ws.Range("K1150").Value = 0
$ws.Range("L1150").Value = 0
$ws.Range("M1150").Value = 0
$ws.Range("N1150").Value = 33
$ws.Range("O1150").Value = 0
$ws.Range("P1150").Value = 0
$ws.Range("Q1150").Value = 0
$ws.Range("R1150").Value = "'"
$ws.Range("R1150").Style = "Normal"
$ws.Range("A1150").NumberFormat = $dateFmt

# Row 1151
$ws.Range("A1151").Value = 45523
$ws.Range("B1151").Value = 312.8337017828379
$ws.Range("C1151").Value = 327.3333016412976
$ws.Range("D1151").Value = 310.1026139101305
$ws.Range("E1151").Value = 321.3248901367188
$ws.Range("F1151").Value = "'"
$ws.Range("F1151").Style = "Normal"
$ws.Range("G1151").Value = 23151861
$ws.Range("H1151").Value = 2024
$ws.Range("I1151").Value = 8
$ws.Range("J1151").Value = 19
$ws.Range("K1151").Value = 0
$ws.Range("L1151").Value = 0
$ws.Range("M1151").Value = 0
$ws.Range("N1151").Value = 34
$ws.Range("O1151").Value = 0
$ws.Range("P1151").Value = 0
$ws.Range("Q1151").Value = 0
$ws.Range("R1151").Value = "'"
$ws.Range("R1151").Style = "Normal"
$ws.Range("A1151").NumberFormat = $dateFmt

# Row 1152
$ws.Range("A1152").Value = 45530
$ws.Range("B1152").Value = 322.7649280749295
$ws.Range("C1152").Value = 325.8436149888774
$ws.Range("D1152").Value = 314.8199452300082
$ws.Range("E1152").Value = 316.8558349609375
$ws.Range("F1152").Value = "'"
$ws.Range("F1152").Style = "Normal"
$ws.Range("G1152").Value = 12545664
$ws.Range("H1152").Value = 2024
$ws.Range("I1152").Value = 8
$ws.Range("J1152").Value = 26
$ws.Range("K1152").Value = 0
$ws.Range("L1152").Value = 0
$ws.Range("M1152").Value = 0
$ws.Range("N1152").Value = 35
$ws.Range("O1152").Value = 0
$ws.Range("P1152").Value = 0
$ws.Range("Q1152").Value = 0
$ws.Range("R1152").Value = "'"
$ws.Range("R1152").Style = "Normal"
$ws.Range("A1152").NumberFormat = $dateFmt

# Row 1153
$ws.Range("A1153").Value = 45537
$ws.Range("B1153").Value = 322.5
$ws.Range("C1153").Value = 330.7999877929688
$ws.Range("D1153").Value = 314.2999877929688
$ws.Range("E1153").Value = 316.1499938964844
$ws.Range("F1153").Value = "'"
$ws.Range("F1153").Style = "Normal"
$ws.Range("G1153").Value = 18081268
$ws.Range("H1153").Value = 2024
$ws.Range("I1153").Value = 9
$ws.Range("J1153").Value = 2
$ws.Range("K1153").Value = 0
$ws.Range("L1153").Value = 0
$ws.Range("M1153").Value = 0
$ws.Range("N1153").Value = 36
$ws.Range("O1153").Value = 0
$ws.Range("P1153").Value = 0
$ws.Range("Q1153").Value = 1
$ws.Range("R1153").Value = "'"
$ws.Range("R1153").Style = "Normal"
$ws.Range("A1153").NumberFormat = $dateFmt

# Row 1154
$ws.Range("A1154").Value = 45544
$ws.Range("B1154").Value = 315.9500122070312
$ws.Range("C1154").Value = 322
$ws.Range("D1154").Value = 308.1000061035156
$ws.Range("E1154").Value = 316.7000122070312
$ws.Range("F1154").Value = "'"
$ws.Range("F1154").Style = "Normal"
$ws.Range("G1154").Value = 11040106
$ws.Range("H1154").Value = 2024
$ws.Range("I1154").Value = 9
$ws.Range("J1154").Value = 9
$ws.Range("K1154").Value = 0
$ws.Range("L1154").Value = 0
$ws.Range("M1154").Value = 0
$ws.Range("N1154").Value = 37
$ws.Range("O1154").Value = 0
$ws.Range("P1154").Value = 0
$ws.Range("Q1154").Value = 2
$ws.Range("R1154").Value = "'"
$ws.Range("R1154").Style = "Normal"
$ws.Range("A1154").NumberFormat = $dateFmt

# Row 1155
$ws.Range("A1155").Value = 45551
$ws.Range("B1155").Value = 317.0499877929688
$ws.Range("C1155").Value = 317.3999938964844
$ws.Range("D1155").Value = 299.1000061035156
$ws.Range("E1155").Value = 314.9500122070312
$ws.Range("F1155").Value = "'"
$ws.Range("F1155").Style = "Normal"
$ws.Range("G1155").Value = 14332426
$ws.Range("H1155").Value = 2024
$ws.Range("I1155").Value = 9
$ws.Range("J1155").Value = 16
$ws.Range("K1155").Value = 0
$ws.Range("L1155").Value = 0
$ws.Range("M1155").Value = 0
$ws.Range("N1155").Value = 38
$ws.Range("O1155").Value = 0
$ws.Range("P1155").Value = 0
$ws.Range("Q1155").Value = 0
$ws.Range("R1155").Value = "'"
$ws.Range("R1155").Style = "Normal"
$ws.Range("A1155").NumberFormat = $dateFmt

# Row 1156
$ws.Range("A1156").Value = 45558
$ws.Range("B1156").Value = 316
$ws.Range("C1156").Value = 317.2000122070312
$ws.Range("D1156").Value = 300
$ws.Range("E1156").Value = 302.3999938964844
$ws.Range("F1156").Value = "'"
$ws.Range("F1156").Style = "Normal"
$ws.Range("G1156").Value = 12319432
$ws.Range("H1156").Value = 2024
$ws.Range("I1156").Value = 9
$ws.Range("J1156").Value = 23
$ws.Range("K1156").Value = 0
$ws.Range("L1156").Value = 0
$ws.Range("M1156").Value = 0
$ws.Range("N1156").Value = 39
$ws.Range("O1156").Value = 0
$ws.Range("P1156").Value = 0
$ws.Range("Q1156").Value = 0
$ws.Range("R1156").Value = "'"
$ws.Range("R1156").Style = "Normal"
$ws.Range("A1156").NumberFormat = $dateFmt

# Row 1157
$ws.Range("A1157").Value = 45565
$ws.Range("B1157").Value = 297.2999877929688
$ws.Range("C1157").Value = 310.4500122070312
$ws.Range("D1157").Value = 291.6499938964844
$ws.Range("E1157").Value = 300.5
$ws.Range("F1157").Value = "'"
$ws.Range("F1157").Style = "Normal"
$ws.Range("G1157").Value = 16573814
$ws.Range("H1157").Value = 2024
$ws.Range("I1157").Value = 9
$ws.Range("J1157").Value = 30
$ws.Range("K1157").Value = 0
$ws.Range("L1157").Value = 0
$ws.Range("M1157").Value = 0
$ws.Range("N1157").Value = 40
$ws.Range("O1157").Value = 0
$ws.Range("P1157").Value = 0
$ws.Range("Q1157").Value = 0
$ws.Range("R1157").Value = "'"
$ws.Range("R1157").Style = "Normal"
$ws.Range("A1157").NumberFormat = $dateFmt

# Row 1158
$ws.Range("A1158").Value = 45572
$ws.Range("B1158").Value = 302.0499877929688
$ws.Range("C1158").Value = 311.9500122070312
$ws.Range("D1158").Value = 289.1499938964844
$ws.Range("E1158").Value = 298.4500122070312
$ws.Range("F1158").Value = "'"
$ws.Range("F1158").Style = "Normal"
$ws.Range("G1158").Value = 15004906
$ws.Range("H1158").Value = 2024
$ws.Range("I1158").Value = 10
$ws.Range("J1158").Value = 7
$ws.Range("K1158").Value = 0
$ws.Range("L1158").Value = 0
$ws.Range("M1158").Value = 0
$ws.Range("N1158").Value = 41
$ws.Range("O1158").Value = 0
$ws.Range("P1158").Value = 0
$ws.Range("Q1158").Value = 0
$ws.Range("R1158").Value = "'"
$ws.Range("R1158").Style = "Normal"
$ws.Range("A1158").NumberFormat = $dateFmt

# Row 1159
$ws.Range("A1159").Value = 45579
$ws.Range("B1159").Value = 299
$ws.Range("C1159").Value = 312.8999938964844
$ws.Range("D1159").Value = 295.2000122070312
$ws.Range("E1159").Value = 310.8999938964844
$ws.Range("F1159").Value = "'"
$ws.Range("F1159").Style = "Normal"
$ws.Range("G1159").Value = 11349734
$ws.Range("H1159").Value = 2024
$ws.Range("I1159").Value = 10
$ws.Range("J1159").Value = 14
$ws.Range("K1159").Value = 0
$ws.Range("L1159").Value = 0
$ws.Range("M1159").Value = 0
$ws.Range("N1159").Value = 42
$ws.Range("O1159").Value = 0
$ws.Range("P1159").Value = 0
$ws.Range("Q1159").Value = 0
$ws.Range("R1159").Value = "'"
$ws.Range("R1159").Style = "Normal"
$ws.Range("A1159").NumberFormat = $dateFmt

# Row 1160
$ws.Range("A1160").Value = 45586
$ws.Range("B1160").Value = 311
$ws.Range("C1160").Value = 311.8999938964844
$ws.Range("D1160").Value = 272.2999877929688
$ws.Range("E1160").Value = 277.6000061035156
$ws.Range("F1160").Value = "'"
$ws.Range("F1160").Style = "Normal"
$ws.Range("G1160").Value = 12233580
$ws.Range("H1160").Value = 2024
$ws.Range("I1160").Value = 10
$ws.Range("J1160").Value = 21
$ws.Range("K1160").Value = 0
$ws.Range("L1160").Value = 0
$ws.Range("M1160").Value = 0
$ws.Range("N1160").Value = 43
$ws.Range("O1160").Value = 0
$ws.Range("P1160").Value = 0
$ws.Range("Q1160").Value = 0
$ws.Range("R1160").Value = "'"
$ws.Range("R1160").Style = "Normal"
$ws.Range("A1160").NumberFormat = $dateFmt

# Row 1161
$ws.Range("A1161").Value = 45593
$ws.Range("B1161").Value = 278
$ws.Range("C1161").Value = 317.9500122070312
$ws.Range("D1161").Value = 270.6499938964844
$ws.Range("E1161").Value = 312.6000061035156
$ws.Range("F1161").Value = "'"
$ws.Range("F1161").Style = "Normal"
$ws.Range("G1161").Value = 15995976
$ws.Range("H1161").Value = 2024
$ws.Range("I1161").Value = 10
$ws.Range("J1161").Value = 28
$ws.Range("K1161").Value = 0
$ws.Range("L1161").Value = 0
$ws.Range("M1161").Value = 0
$ws.Range("N1161").Value = 44
$ws.Range("O1161").Value = 0
$ws.Range("P1161").Value = 0
$ws.Range("Q1161").Value = 0
$ws.Range("R1161").Value = "'"
$ws.Range("R1161").Style = "Normal"
$ws.Range("A1161").NumberFormat = $dateFmt

# Row 1162
$ws.Range("A1162").Value = 45600
$ws.Range("B1162").Value = 317.7000122070312
$ws.Range("C1162").Value = 321.8999938964844
$ws.Range("D1162").Value = 301.0499877929688
$ws.Range("E1162").Value = 308.1000061035156
$ws.Range("F1162").Value = "'"
$ws.Range("F1162").Style = "Normal"
$ws.Range("G1162").Value = 19774990
$ws.Range("H1162").Value = 2024
$ws.Range("I1162").Value = 11
$ws.Range("J1162").Value = 4
$ws.Range("K1162").Value = 0
$ws.Range("L1162").Value = 0
$ws.Range("M1162").Value = 0
$ws.Range("N1162").Value = 45
$ws.Range("O1162").Value = 0
$ws.Range("P1162").Value = 0
$ws.Range("Q1162").Value = 0
$ws.Range("R1162").Value = "'"
$ws.Range("R1162").Style = "Normal"
$ws.Range("A1162").NumberFormat = $dateFmt

# Row 1163
$ws.Range("A1163").Value = 45607
$ws.Range("B1163").Value = 306.25
$ws.Range("C1163").Value = 309.5
$ws.Range("D1163").Value = 277.3500061035156
$ws.Range("E1163").Value = 280
$ws.Range("F1163").Value = "'"
$ws.Range("F1163").Style = "Normal"
$ws.Range("G1163").Value = 9640772
$ws.Range("H1163").Value = 2024
$ws.Range("I1163").Value = 11
$ws.Range("J1163").Value = 11
$ws.Range("K1163").Value = 0
$ws.Range("L1163").Value = 0
$ws.Range("M1163").Value = 0
$ws.Range("N1163").Value = 46
$ws.Range("O1163").Value = 0
$ws.Range("P1163").Value = 0
$ws.Range("Q1163").Value = 0
$ws.Range("R1163").Value = "'"
$ws.Range("R1163").Style = "Normal"
$ws.Range("A1163").NumberFormat = $dateFmt

# Row 1164
$ws.Range("A1164").Value = 45614
$ws.Range("B1164").Value = 279.9500122070312
$ws.Range("C1164").Value = 288.2000122070312
$ws.Range("D1164").Value = 269.25
$ws.Range("E1164").Value = 283.7999877929688
$ws.Range("F1164").Value = "'"
$ws.Range("F1164").Style = "Normal"
$ws.Range("G1164").Value = 9457688
$ws.Range("H1164").Value = 2024
$ws.Range("I1164").Value = 11
$ws.Range("J1164").Value = 18
$ws.Range("K1164").Value = 0
$ws.Range("L1164").Value = 0
$ws.Range("M1164").Value = 0
$ws.Range("N1164").Value = 47
$ws.Range("O1164").Value = 2
$ws.Range("P1164").Value = 0
$ws.Range("Q1164").Value = 0
$ws.Range("R1164").Value = "'"
$ws.Range("R1164").Style = "Normal"
$ws.Range("A1164").NumberFormat = $dateFmt

# Row 1165
$ws.Range("A1165").Value = 45621
$ws.Range("B1165").Value = 290
$ws.Range("C1165").Value = 318.25
$ws.Range("D1165").Value = 289.1499938964844
$ws.Range("E1165").Value = 309.7000122070312
$ws.Range("F1165").Value = "'"
$ws.Range("F1165").Style = "Normal"
$ws.Range("G1165").Value = 20543851
$ws.Range("H1165").Value = 2024
$ws.Range("I1165").Value = 11
$ws.Range("J1165").Value = 25
$ws.Range("K1165").Value = 0
$ws.Range("L1165").Value = 0
$ws.Range("M1165").Value = 0
$ws.Range("N1165").Value = 48
$ws.Range("O1165").Value = 0
$ws.Range("P1165").Value = 0
$ws.Range("Q1165").Value = 0
$ws.Range("R1165").Value = "'"
$ws.Range("R1165").Style = "Normal"
$ws.Range("A1165").NumberFormat = $dateFmt

# Row 1166
$ws.Range("A1166").Value = 45628
$ws.Range("B1166").Value = 311.8999938964844
$ws.Range("C1166").Value = 326.4500122070312
$ws.Range("D1166").Value = 304.1499938964844
$ws.Range("E1166").Value = 314.1000061035156
$ws.Range("F1166").Value = "'"
$ws.Range("F1166").Style = "Normal"
$ws.Range("G1166").Value = 20385857
$ws.Range("H1166").Value = 2024
$ws.Range("I1166").Value = 12
$ws.Range("J1166").Value = 2
$ws.Range("K1166").Value = 0
$ws.Range("L1166").Value = 0
$ws.Range("M1166").Value = 0
$ws.Range("N1166").Value = 49
$ws.Range("O1166").Value = 0
$ws.Range("P1166").Value = 0
$ws.Range("Q1166").Value = 0
$ws.Range("R1166").Value = "'"
$ws.Range("R1166").Style = "Normal"
$ws.Range("A1166").NumberFormat = $dateFmt

# Row 1167
$ws.Range("A1167").Value = 45635
$ws.Range("B1167").Value = 315.7000122070312
$ws.Range("C1167").Value = 317.25
$ws.Range("D1167").Value = 294.1000061035156
$ws.Range("E1167").Value = 302.7000122070312
$ws.Range("F1167").Value = "'"
$ws.Range("F1167").Style = "Normal"
$ws.Range("G1167").Value = 10590298
$ws.Range("H1167").Value = 2024
$ws.Range("I1167").Value = 12
$ws.Range("J1167").Value = 9
$ws.Range("K1167").Value = 0
$ws.Range("L1167").Value = 0
$ws.Range("M1167").Value = 0
$ws.Range("N1167").Value = 50
$ws.Range("O1167").Value = 0
$ws.Range("P1167").Value = 0
$ws.Range("Q1167").Value = 0
$ws.Range("R1167").Value = "'"
$ws.Range("R1167").Style = "Normal"
$ws.Range("A1167").NumberFormat = $dateFmt

# Row 1168
$ws.Range("A1168").Value = 45642
$ws.Range("B1168").Value = 302.7000122070312
$ws.Range("C1168").Value = 315.3999938964844
$ws.Range("D1168").Value = 286.5
$ws.Range("E1168").Value = 287.5499877929688
$ws.Range("F1168").Value = "'"
$ws.Range("F1168").Style = "Normal"
$ws.Range("G1168").Value = 16152117
$ws.Range("H1168").Value = 2024
$ws.Range("I1168").Value = 12
$ws.Range("J1168").Value = 16
$ws.Range("K1168").Value = 0
$ws.Range("L1168").Value = 0
$ws.Range("M1168").Value = 0
$ws.Range("N1168").Value = 51
$ws.Range("O1168").Value = 0
$ws.Range("P1168").Value = 0
$ws.Range("Q1168").Value = 0
$ws.Range("R1168").Value = "'"
$ws.Range("R1168").Style = "Normal"
$ws.Range("A1168").NumberFormat = $dateFmt

# Row 1169
$ws.Range("A1169").Value = 45649
$ws.Range("B1169").Value = 289.9500122070312
$ws.Range("C1169").Value = 291.9500122070312
$ws.Range("D1169").Value = 278.6499938964844
$ws.Range("E1169").Value = 279.1499938964844
$ws.Range("F1169").Value = "'"
$ws.Range("F1169").Style = "Normal"
$ws.Range("G1169").Value = 6620700
$ws.Range("H1169").Value = 2024
$ws.Range("I1169").Value = 12
$ws.Range("J1169").Value = 23
$ws.Range("K1169").Value = 0
$ws.Range("L1169").Value = 0
$ws.Range("M1169").Value = 0
$ws.Range("N1169").Value = 52
$ws.Range("O1169").Value = 0
$ws.Range("P1169").Value = 0
$ws.Range("Q1169").Value = 0
$ws.Range("R1169").Value = "'"
$ws.Range("R1169").Style = "Normal"
$ws.Range("A1169").NumberFormat = $dateFmt
